{"js": "// \"Mensajes edicion con graficos.\"\n//\n// Updates the quotation's date, the work's title (heading + \"T\u00edtulo:\"\n// value), the authors list, the page count, the cover/flap placeholder\n// text, the imprenta name and the table's TOTAL($) value.\n//\n// Each replacement is scoped to the single paragraph that owns the text\n// (found via a short, unique prefix) and then to the specific sub-range\n// inside that paragraph, so only the intended run's text changes and\n// sibling runs (e.g. the bold \"Label: \" run sharing the paragraph) keep\n// their original formatting untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfunction findParagraph(startsWith) {\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.indexOf(startsWith) === 0) {\n      return paragraphs.items[i];\n    }\n  }\n  throw new Error(\"Paragraph not found for: \" + startsWith);\n}\n\nasync function replaceInParagraph(startsWith, oldText, newText) {\n  const para = findParagraph(startsWith);\n  const hits = para.search(oldText, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  if (hits.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText + \" (in paragraph starting '\" + startsWith + \"')\");\n  }\n  hits.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1) Quotation date heading.\nawait replaceInParagraph(\n  \"Guayaquil, 19 de enero de 2018\",\n  \"Guayaquil, 19 de enero de 2018\",\n  \"Guayaquil, 02 de febrero de 2018\"\n);\n\n// 2) Work title heading (bold, centered).\nawait replaceInParagraph(\n  \"Gesti\u00f3n para la formaci\u00f3n en gobernabilidad en el escenario ecuatoriano\",\n  \"Gesti\u00f3n para la formaci\u00f3n en gobernabilidad en el escenario ecuatoriano\",\n  \"Student perceptions of the use of SIDWEB for Learning English writing skills in an ecuadorian university\"\n);\n\n// 3) \"Cotizaci\u00f3n solicitada...\" sentence (month reference).\nawait replaceInParagraph(\n  \"Cotizaci\u00f3n solicitada, en mes de enero de 2018\",\n  \"Cotizaci\u00f3n solicitada, en mes de enero de 2018, de acuerdo con las siguientes caracter\u00edsticas:\",\n  \"Cotizaci\u00f3n solicitada, en mes de febrero de 2018, de acuerdo con las siguientes caracter\u00edsticas:\"\n);\n\n// 4) \"T\u00edtulo: \" value run.\nawait replaceInParagraph(\n  \"T\u00edtulo: \",\n  \"Gesti\u00f3n para la formaci\u00f3n en gobernabilidad en el escenario ecuatoriano\",\n  \"Student perceptions of the use of SIDWEB for Learning English writing skills in an ecuadorian university\"\n);\n\n// 5) \"Autores: \" value run.\nawait replaceInParagraph(\n  \"Autores: \",\n  \"Mauro Toscanini Segale, Uriel Castillo Nazareno, Jack Ch\u00e1vez Garc\u00eda, Teresa Alc\u00edvar Avil\u00e9s,  Tamara Proenza D\u00edaz.\",\n  \"Roxana Fern\u00e1ndez Berducci.\"\n);\n\n// 6) \"N\u00famero de p\u00e1ginas: \" value run.\nawait replaceInParagraph(\"N\u00famero de p\u00e1ginas: \", \"1\", \"500\");\n\n// 7) \"Cubierta: \" value run.\nawait replaceInParagraph(\"Cubierta: \", \"adasd\", \"fg\");\n\n// 8) \"Solapas: \" value run.\nawait replaceInParagraph(\"Solapas: \", \"asdasd\", \"fgsd\");\n\n// 9) Table cell: imprenta name.\nawait replaceInParagraph(\"imprentaui\", \"imprentaui\", \"sdfsdf\");\n\n// 10) Table cell: TOTAL ($) value.\nawait replaceInParagraph(\"$67.2\", \"$67.2\", \"$60\");\n", "ps1": "# \"Mensajes edicion con graficos.\"\n#\n# Updates the quotation's date, the work's title (heading + \"T\u00edtulo:\"\n# value), the authors list, the page count, the cover/flap placeholder\n# text, the imprenta name and the table's TOTAL($) value.\n#\n# Each replacement is scoped to the single paragraph that owns the text\n# (found via a short, unique prefix) and run through that paragraph's own\n# Find/Replace, so only the intended text changes and sibling runs (e.g.\n# the bold \"Label: \" run sharing the paragraph) keep their original\n# formatting untouched.\n#\n# Literal strings use single quotes throughout so none of them (e.g. the\n# '$67.2' / '$60' amounts) are misread as PowerShell variable expansion.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphByPrefix($prefix) {\n  $count = $d.Paragraphs.Count\n  for ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text.StartsWith($prefix)) {\n      return $p\n    }\n  }\n  throw \"paragraph not found for prefix: $prefix\"\n}\n\nfunction Replace-InParagraph($prefix, $oldText, $newText) {\n  $p = Get-ParagraphByPrefix($prefix)\n  $find = $p.Range.Find\n  $ok = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $ok) {\n    throw \"Find/Replace failed for '$oldText' in paragraph starting '$prefix'\"\n  }\n}\n\n# 1) Quotation date heading.\nReplace-InParagraph 'Guayaquil, 19 de enero de 2018' 'Guayaquil, 19 de enero de 2018' 'Guayaquil, 02 de febrero de 2018'\n\n# 2) Work title heading (bold, centered).\nReplace-InParagraph 'Gesti\u00f3n para la formaci\u00f3n en gobernabilidad en el escenario ecuatoriano' 'Gesti\u00f3n para la formaci\u00f3n en gobernabilidad en el escenario ecuatoriano' 'Student perceptions of the use of SIDWEB for Learning English writing skills in an ecuadorian university'\n\n# 3) \"Cotizaci\u00f3n solicitada...\" sentence (month reference).\nReplace-InParagraph 'Cotizaci\u00f3n solicitada, en mes de enero de 2018' 'Cotizaci\u00f3n solicitada, en mes de enero de 2018, de acuerdo con las siguientes caracter\u00edsticas:' 'Cotizaci\u00f3n solicitada, en mes de febrero de 2018, de acuerdo con las siguientes caracter\u00edsticas:'\n\n# 4) \"T\u00edtulo: \" value run.\nReplace-InParagraph 'T\u00edtulo: ' 'Gesti\u00f3n para la formaci\u00f3n en gobernabilidad en el escenario ecuatoriano' 'Student perceptions of the use of SIDWEB for Learning English writing skills in an ecuadorian university'\n\n# 5) \"Autores: \" value run.\nReplace-InParagraph 'Autores: ' 'Mauro Toscanini Segale, Uriel Castillo Nazareno, Jack Ch\u00e1vez Garc\u00eda, Teresa Alc\u00edvar Avil\u00e9s,  Tamara Proenza D\u00edaz.' 'Roxana Fern\u00e1ndez Berducci.'\n\n# 6) \"N\u00famero de p\u00e1ginas: \" value run.\nReplace-InParagraph 'N\u00famero de p\u00e1ginas: ' '1' '500'\n\n# 7) \"Cubierta: \" value run.\nReplace-InParagraph 'Cubierta: ' 'adasd' 'fg'\n\n# 8) \"Solapas: \" value run.\nReplace-InParagraph 'Solapas: ' 'asdasd' 'fgsd'\n\n# 9) Table cell: imprenta name.\nReplace-InParagraph 'imprentaui' 'imprentaui' 'sdfsdf'\n\n# 10) Table cell: TOTAL ($) value.\nReplace-InParagraph '$67.2' '$67.2' '$60'\n"}
